# eFP-Seq Browser bulk_template.xlsx submission fix
# -----------------------------------------------------------------
# Commit: "Excel conversion and generating data bug fixes" -- adds a
# missing RNA-Seq submission entry (row 3, "Arabidopsis thaliana
# Flowers 12-14" / SRR3581866) to the "FILL Data" worksheet's Table1,
# which also grows the table/autofilter range from A1:S3 to A1:S4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FILL Data")
$wsBasic = $wb.Worksheets.Item("FILL Basic")

# Grow the worksheet table by one row -- this both extends the
# table/autoFilter ref to A1:S4 and gives us a blank row 4 to fill in.
$tbl = $ws.ListObjects.Item(1)
[void]$tbl.ListRows.Add()
$r = 4

# --- formatting -----------------------------------------------------
# A handful of columns in the new row reuse formats that already exist
# elsewhere in the workbook (plain "Arial, no fill/border" cells from
# the FILL Basic sheet, and the heavily-used "Amazon AWS"/"Species"
# style from this same table). Copy those formats across first so the
# saved file reuses the existing style records instead of minting
# duplicates.
$wsBasic.Cells.Item(2, 4).Copy()
$ws.Cells.Item($r, 4).PasteSpecial(-4122)   # Record Number*      -> style like D2 (FILL Basic)
$ws.Cells.Item($r, 11).PasteSpecial(-4122)  # Tissue*             -> style like D2 (FILL Basic)
$ws.Cells.Item($r, 12).PasteSpecial(-4122)  # Tissue subunit*     -> style like D2 (FILL Basic)

$ws.Cells.Item(2, 6).Copy()
$ws.Cells.Item($r, 6).PasteSpecial(-4122)   # Repository type*    -> existing "Amazon AWS" style
$ws.Cells.Item($r, 10).PasteSpecial(-4122)  # Species*            -> existing "Arabidopsis thaliana" style

$excel.CutCopyMode = 0

# --- values -----------------------------------------------------------
$ws.Cells.Item($r, 1).Value2 = 3
$ws.Cells.Item($r, 2).Value2 = "Arabidopsis thaliana Flowers 12-14"
$ws.Cells.Item($r, 3).Value2 = "1 ILLUMINA (Illumina HiSeq 2000) run: 25.1M spots, 1.3G bases, 806.8Mb downloads"
$ws.Cells.Item($r, 4).Value2 = "SRR3581866"
$ws.Cells.Item($r, 5).Value2 = "http://bar.utoronto.ca/~asullivan/data/SRR3581866/"
$ws.Cells.Item($r, 6).Value2 = "Amazon AWS"
$ws.Cells.Item($r, 7).Value2 = "https://www.ncbi.nlm.nih.gov/pubmed/27549386"
$ws.Cells.Item($r, 8).Value2 = "https://trace.ncbi.nlm.nih.gov/Traces/sra/?run=SRR3581866"
$ws.Cells.Item($r, 9).Value2 = 25081651
$ws.Cells.Item($r, 10).Value2 = "Arabidopsis thaliana"
$ws.Cells.Item($r, 11).Value2 = "Flower"
$ws.Cells.Item($r, 12).Value2 = "flower"
# M4:S4 (Controls / Replicate Control 1-6) are intentionally left blank.

# The remaining new-row cells (Title*, Description*, RNA-Seq link,
# Publication Link, SRA/NCBI Link, Total Reads Mapped*, and the blank
# Controls/Replicate columns) get a plain black-Calibri, no-fill,
# no-border look that isn't used anywhere else yet in the workbook, so
# it becomes a brand-new cell style on save.
$newStyleCols = @(2, 3, 5, 7, 8, 9, 13, 14, 15, 16, 17, 18, 19)
foreach ($col in $newStyleCols) {
    $cell = $ws.Cells.Item($r, $col)
    $cell.Font.Color = 0
    $cell.Font.ThemeFont = 0
}

# --- selection ---------------------------------------------------------
# Reflect where the submitter's view ended up after adding the row.
$ws.Range("M4").Select()
